$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 21:13"

# Country name realignment (Tunez / Namibia moved earlier in the country ordering,
# which shifts which country each of these rows now represents)
$ws.Range("A79").Value = "Tunez"
$ws.Range("A80").Value = "Bosnia y Herzegovina"
$ws.Range("A81").Value = "El Salvador"
$ws.Range("A82").Value = "Australia"
$ws.Range("A101").Value = "Namibia"
$ws.Range("A102").Value = "Finlandia"

# Updated statistics per country row
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7926376
$ws.Range("C4").Value = 31898
$ws.Range("D4").Value = 5078346
$ws.Range("E4").Value = 2629001
$ws.Range("G4").Value = 382
$ws.Range("H4").Value = 219029

# Row 5 - India
$ws.Range("B5").Value = 7037694
$ws.Range("C5").Value = 60686
$ws.Range("D5").Value = 6046028
$ws.Range("E5").Value = 883619
$ws.Range("G5").Value = 597
$ws.Range("H5").Value = 108047

# Row 13 - Francia
$ws.Range("B13").Value = 718873
$ws.Range("C13").Value = 26896
$ws.Range("E13").Value = 585462

# Row 25 - Alemania
$ws.Range("B25").Value = 322755
$ws.Range("C25").Value = 2277
$ws.Range("E25").Value = 39565

# Row 29 - Canada
$ws.Range("B29").Value = 180142
$ws.Range("C29").Value = 2025
$ws.Range("D29").Value = 151345
$ws.Range("E29").Value = 19189
$ws.Range("G29").Value = 23
$ws.Range("H29").Value = 9608

# Row 59 - Uzbekistan
$ws.Range("B59").Value = 60776
$ws.Range("C59").Value = 434
$ws.Range("D59").Value = 57704
$ws.Range("E59").Value = 2570
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 502

# Row 79 - Tunez (after realignment)
$ws.Range("B79").Value = 31259
$ws.Range("C79").Value = 4360
$ws.Range("D79").Value = 5032
$ws.Range("E79").Value = 25771
$ws.Range("G79").Value = 47
$ws.Range("H79").Value = 456

# Row 80 - Bosnia y Herzegovina (after realignment)
$ws.Range("B80").Value = 30345
$ws.Range("C80").Value = 428
$ws.Range("D80").Value = 23370
$ws.Range("E80").Value = 6048
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 927

# Row 81 - El Salvador (after realignment)
$ws.Range("B81").Value = 29951
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 24995
$ws.Range("E81").Value = 4069
$ws.Range("G81").Value = 6
$ws.Range("H81").Value = 887

# Row 82 - Australia (after realignment)
$ws.Range("B82").Value = 27244
$ws.Range("C82").Value = 15
$ws.Range("D82").Value = 24987
$ws.Range("E82").Value = 1360
$ws.Range("H82").Value = 897

# Row 101 - Namibia (after realignment)
$ws.Range("B101").Value = 11891
$ws.Range("C101").Value = 62
$ws.Range("D101").Value = 9807
$ws.Range("E101").Value = 1956
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 128

# Row 102 - Finlandia (after realignment)
$ws.Range("B102").Value = 11849
$ws.Range("C102").Value = 269
$ws.Range("D102").Value = 8500
$ws.Range("E102").Value = 3003
$ws.Range("H102").Value = 346

# Row 138 - Aruba
$ws.Range("B138").Value = 4188
$ws.Range("C138").Value = 21
$ws.Range("D138").Value = 3753
$ws.Range("E138").Value = 403

# Row 147 - Mali
$ws.Range("B147").Value = 3273
$ws.Range("C147").Value = 25
$ws.Range("D147").Value = 2525
$ws.Range("E147").Value = 616
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 132

# Row 183 - Eritrea
$ws.Range("B183").Value = 414
$ws.Range("C183").Value = 9
$ws.Range("D183").Value = 372
$ws.Range("E183").Value = 42
